$d = $word.ActiveDocument

# Append a new note paragraph right after the last picture paragraph
# (and before the trailing blank paragraph), before touching proofing
# flags so the new run does not inherit <w:noProof/>.
$lastPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newPara.Range.Text = "Grab filename with score and sort desc, grab data with condition score === 8.0, group by file name"

# Mark every inline picture's run as "do not spell/grammar check"
# (adds <w:noProof/> to that run's rPr) -- mirrors what Word does when
# a picture run is (re)inserted / the doc is reflowed.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes($i)
    $shape.Range.NoProofing = 1
}
